$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.780.38'
$ws.Range("E2").Value = '  +2.91%  '
# Row 3
$ws.Range("D3").Value = '2.625.92'
$ws.Range("E3").Value = '  +2.07%  '
# Row 4
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.58'
$ws.Range("E5").Value = '  -0.71%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.46'
$ws.Range("E6").Value = '  -0.15%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.12%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +0.74%  '
# Row 9
$ws.Range("D9").Value = '2.627.38'
$ws.Range("E9").Value = '  +2.01%  '
# Row 10
$ws.Range("E10").Value = '  -2.53%  '
# Row 11
$ws.Range("E11").Value = '  +2.71%  '
# Row 12
$ws.Range("E12").Value = '  -3.71%  '
# Row 13
$ws.Range("E13").Value = '  +7.09%  '
# Row 14
$ws.Range("D14").Value = '3.071.29'
$ws.Range("E14").Value = '  +1.50%  '
# Row 15
$ws.Range("D15").Value = '60.805.31'
$ws.Range("E15").Value = '  +2.88%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.44'
$ws.Range("E16").Value = '  +5.22%  '
# Row 17
$ws.Range("E17").Value = '  +2.94%  '
# Row 18
$ws.Range("D18").Value = '2.616.38'
$ws.Range("E18").Value = '  +1.28%  '
# Row 19
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.19'
$ws.Range("E19").Value = '  +9.14%  '
# Row 20
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.65'
$ws.Range("E20").Value = '  +2.80%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.80'
$ws.Range("E21").Value = '  +3.91%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.06'
$ws.Range("E22").Value = '  +12.43%  '
# Row 23
$ws.Range("E23").Value = '  +0.17%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.525'
$ws.Range("E24").Value = '  +13.37%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.30'
$ws.Range("E25").Value = '  -0.45%  '
# Row 26
$ws.Range("E26").Value = '  -0.87%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("E27").Value = '  -0.95%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("E28").Value = '  +5.95%  '
# Row 29
$ws.Range("D29").Value = '0.0₃0792'
$ws.Range("E29").Value = '  +1.69%  '
# Row 30
$ws.Range("E30").Value = '  +8.33%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  +0.00%  '
# Row 32
$ws.Range("E32").Value = '  +4.91%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '161.94'
$ws.Range("E33").Value = '  +1.23%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.49'
$ws.Range("E34").Value = '  +2.95%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.24'
$ws.Range("E35").Value = '  +5.64%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.957'
$ws.Range("E36").Value = '  +9.43%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.20'
$ws.Range("E37").Value = '  +3.58%  '
# Row 38
$ws.Range("E38").Value = '  +7.30%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.75'
$ws.Range("E39").Value = '  +0.98%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.848'
$ws.Range("E40").Value = '  -2.80%  '
# Row 41
$ws.Range("E41").Value = '  +3.55%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '296.81'
$ws.Range("E42").Value = '  +0.64%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '140.39'
$ws.Range("E43").Value = '  +6.95%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  -0.09%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0985'
$ws.Range("E45").Value = '  +0.96%  '
# Row 46
$ws.Range("E46").Value = '  +2.08%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0547'
$ws.Range("E47").Value = '  +2.10%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.51'
$ws.Range("E48").Value = '  +2.12%  '
# Row 49
$ws.Range("E49").Value = '  +3.47%  '
# Row 50
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.70'
$ws.Range("E50").Value = '  +0.54%  '
# Row 51
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.68'
$ws.Range("E51").Value = '  +6.43%  '
